$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 134 (shifts existing rows 134-206 down to 135-207,
# just like Excel's Range.EntireRow.Insert / Rows(n).Insert()).
$ws.Rows(134).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(134, 1).Value = 4
$ws.Cells.Item(134, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(134, 3).Value = "Los Lagos"
$ws.Cells.Item(134, 4).Value = 44704
$ws.Cells.Item(134, 5).Value = 10
$ws.Cells.Item(134, 6).Value = "Fruta"
$ws.Cells.Item(134, 7).Value = 100102
$ws.Cells.Item(134, 8).Value = "Cítricos"
$ws.Cells.Item(134, 9).Value = 100102004
$ws.Cells.Item(134, 10).Value = "Mandarina"
$ws.Cells.Item(134, 11).Value = "Murcott"
$ws.Cells.Item(134, 12).Value = "Segunda"
$ws.Cells.Item(134, 13).Value = 400
$ws.Cells.Item(134, 14).Value = 11000
$ws.Cells.Item(134, 15).Value = 11000
$ws.Cells.Item(134, 16).Value = 11000
$ws.Cells.Item(134, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(134, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(134, 19).Value = 1100
$ws.Cells.Item(134, 20).Value = 10
